# The deck's design theme is switched from the custom "Integral" theme to
# the built-in "Office Theme" palette (Design tab -> Themes -> Office).
# This rewrites the 12 theme colors (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) on the presentation's one-and-only slide-master theme.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$theme = $m.Theme
$cs = $theme.ThemeColorScheme

# PpColorSchemeIndex order: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2
# 7=accent3 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
# .RGB takes a BGR-packed integer (R + G*256 + B*65536), i.e. the same
# encoding PowerPoint's own RGB() helper produces.

$cs.Item(1).RGB  = 0x00 + 0x00*256 + 0x00*65536   # dk1      000000
$cs.Item(2).RGB  = 0xFF + 0xFF*256 + 0xFF*65536   # lt1      FFFFFF
$cs.Item(3).RGB  = 0x44 + 0x54*256 + 0x6A*65536   # dk2      44546A
$cs.Item(4).RGB  = 0xE7 + 0xE6*256 + 0xE6*65536   # lt2      E7E6E6
$cs.Item(5).RGB  = 0x5B + 0x9B*256 + 0xD5*65536   # accent1  5B9BD5
$cs.Item(6).RGB  = 0xED + 0x7D*256 + 0x31*65536   # accent2  ED7D31
$cs.Item(7).RGB  = 0xA5 + 0xA5*256 + 0xA5*65536   # accent3  A5A5A5
$cs.Item(8).RGB  = 0xFF + 0xC0*256 + 0x00*65536   # accent4  FFC000
$cs.Item(9).RGB  = 0x44 + 0x72*256 + 0xC4*65536   # accent5  4472C4
$cs.Item(10).RGB = 0x70 + 0xAD*256 + 0x47*65536   # accent6  70AD47
$cs.Item(11).RGB = 0x05 + 0x63*256 + 0xC1*65536   # hlink    0563C1
$cs.Item(12).RGB = 0x95 + 0x4F*256 + 0x72*65536   # folHlink 954F72

$theme.Name = "Office Theme"
